$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.919331073760986
$ws.Range("B1").Value = 2.93938684463501
$ws.Range("C1").Value = 1.847727179527283
$ws.Range("D1").Value = 1.57336151599884
$ws.Range("E1").Value = 1.449084401130676
